$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "'10"
$ws.Range("D3").Value = "'23.01"
$ws.Range("G3").Value = "'10"
$ws.Range("D4").Value = "'5.408"
$ws.Range("G4").Value = "'10"
$ws.Range("D5").Value = "'0.05921"
$ws.Range("G5").Value = "'10"
$ws.Range("D6").Value = "'3.452"
$ws.Range("G6").Value = "'10"
$ws.Range("D7").Value = "'6.543"
$ws.Range("G7").Value = "'10"
$ws.Range("D8").Value = "'0.8110"
$ws.Range("G8").Value = "'10"
$ws.Range("D9").Value = "'0.9088"
$ws.Range("G9").Value = "'10"
$ws.Range("D10").Value = "'0.1407"
$ws.Range("G10").Value = "'10"
$ws.Range("D11").Value = "'0.07337"
$ws.Range("G11").Value = "'10"
$ws.Range("G12").Value = "'10"
$ws.Range("D13").Value = "'0.03048"
$ws.Range("G13").Value = "'10"
$ws.Range("D14").Value = "'0.09347"
$ws.Range("G14").Value = "'10"
$ws.Range("D15").Value = "'3.857"
$ws.Range("G15").Value = "'10"
$ws.Range("D16").Value = "'0.001575"
$ws.Range("G16").Value = "'10"
$ws.Range("D17").Value = "'0.04672"
$ws.Range("G17").Value = "'10"
$ws.Range("D18").Value = "'0.01123"
$ws.Range("E18").Value = "'17OneONEBestin24h"
$ws.Range("G18").Value = "'10"
$ws.Range("D19").Value = "'0.006143"
$ws.Range("G19").Value = "'10"
$ws.Range("D20").Value = "'0.004973"
$ws.Range("G20").Value = "'10"
$ws.Range("D21").Value = "'0.0009830"
$ws.Range("G21").Value = "'10"
$ws.Range("D22").Value = "'0.00009408"
$ws.Range("G22").Value = "'10"
$ws.Range("D23").Value = "'3.607"
$ws.Range("G23").Value = "'10"
$ws.Range("D24").Value = "'2.151"
$ws.Range("G24").Value = "'10"
$ws.Range("G25").Value = "'10"
$ws.Range("G26").Value = "'10"
$ws.Range("D27").Value = "'0.0002902"
$ws.Range("G27").Value = "'10"
$ws.Range("G28").Value = "'10"
$ws.Range("G29").Value = "'10"
$ws.Range("G30").Value = "'10"
$ws.Range("G31").Value = "'10"
$ws.Range("G32").Value = "'10"
$ws.Range("G33").Value = "'10"
$ws.Range("G34").Value = "'10"
$ws.Range("G35").Value = "'10"
$ws.Range("G36").Value = "'10"
$ws.Range("G37").Value = "'10"
$ws.Range("G38").Value = "'10"
$ws.Range("G39").Value = "'10"
$ws.Range("D40").Value = "'0.03964"
$ws.Range("G40").Value = "'10"
$ws.Range("D41").Value = "'0.006203"
$ws.Range("E41").Value = "'40KickTokenKICK"
$ws.Range("G41").Value = "'10"
$ws.Range("D42").Value = "'0.1075"
$ws.Range("G42").Value = "'10"
$ws.Range("D43").Value = "'0.003002"
$ws.Range("G43").Value = "'10"
$ws.Range("D44").Value = "'0.008917"
$ws.Range("G44").Value = "'10"
$ws.Range("D45").Value = "'0.00005263"
$ws.Range("G45").Value = "'10"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("G46").Value = "'10"
$ws.Range("D47").Value = "'0.7182"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("G47").Value = "'10"
$ws.Range("D48").Value = "'0.002259"
$ws.Range("G48").Value = "'10"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("G49").Value = "'10"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'49SpecialPowerGoldSPG"
$ws.Range("G50").Value = "'10"
$ws.Range("G51").Value = "'10"
